$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 107, pushing the existing rows 107-112 down to 108-113.
$ws.Rows.Item(107).Insert()

# Populate the newly inserted row 107 with a new weekly data point
# (same series as the row that used to be there, but with a new date).
$ws.Range("A107").Value = 9
$ws.Range("B107").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C107").Value = "Metropolitana"
$ws.Range("D107").Value = 45147
$ws.Range("E107").Value = 13
$ws.Range("F107").Value = 100112029
$ws.Range("G107").Value = "Orégano"
$ws.Range("H107").Value = "Sin especificar"
$ws.Range("I107").Value = "Primera"
$ws.Range("J107").Value = 16
$ws.Range("K107").Value = 21000
$ws.Range("L107").Value = 21000
$ws.Range("M107").Value = 21000
$ws.Range("N107").Value = "$/docena de atados"
$ws.Range("O107").Value = "Región Metropolitana"
$ws.Range("P107").Value = 7000
$ws.Range("Q107").Value = 3
$ws.Range("R107").Value = "Hortaliza"
